$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# General info updates
$ws.Range("A5").Value = "Issue date: 03/12/2020 11:53:13"
$ws.Range("A6").Value = "Python version: Python 3.7.3"

# Row 17: Doc2Vec -> Word2Vec, with updated metrics
$ws.Range("C17").Value = "Word2VecTransfomer"
$ws.Range("E17").Value = "None"
$ws.Range("G17").Value = "70.1*"
$ws.Range("H17").Value = "68.68*"
$ws.Range("I17").Value = "62.65*"
$ws.Range("J17").Value = "72.32*"

# Swap significance colors on I17/J17 (I17 red->blue, J17 blue->red)
$ws.Range("I17").Font.Color = 16711680
$ws.Range("J17").Font.Color = 255

# Column C width: 17.71 -> 18.71 (characters); engine rounds to pixel grid
$ws.Columns("C").ColumnWidth = 17.8
